$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 3.7
$ws.Range("G2").Value = 3.95
$ws.Range("H2").Value = 3.8
$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 2.06
$ws.Range("K2").Value = 2.1
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 0
$ws.Range("V2").Value = 1.31
$ws.Range("W2").Value = 1.34
$ws.Range("AH2").Value = 2.58
$ws.Range("AI2").Value = 6
$ws.Range("AL2").Value = 6.6
$ws.Range("AM2").Value = 12.5
$ws.Range("AN2").Value = 16
$ws.Range("AO2").Value = 17.5
$ws.Range("F3").Value = 2.14
$ws.Range("G3").Value = 3.9
$ws.Range("H3").Value = 3.5
$ws.Range("J3").Value = 2.72
$ws.Range("K3").Value = 3.2
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 1.1
$ws.Range("O3").Value = 1.01
$ws.Range("P3").Value = 1.24
$ws.Range("Q3").Value = 1.03
$ws.Range("R3").Value = 1.21
$ws.Range("S3").Value = 1.03
$ws.Range("T3").Value = 1.9
$ws.Range("U3").Value = 1.76
$ws.Range("V3").Value = 1.02
$ws.Range("W3").Value = 1.03
$ws.Range("X3").Value = 19
$ws.Range("Y3").Value = 1000
$ws.Range("Z3").Value = 1000
$ws.Range("AA3").Value = 1000
$ws.Range("AB3").Value = 1000
$ws.Range("AC3").Value = 960
$ws.Range("AE3").Value = 1000
$ws.Range("AF3").Value = 1000
$ws.Range("AG3").Value = 12
$ws.Range("AH3").Value = 23
$ws.Range("AI3").Value = 970
$ws.Range("AJ3").Value = 1000
$ws.Range("AK3").Value = 1000
$ws.Range("AL3").Value = 1000
$ws.Range("AM3").Value = 1000
$ws.Range("AN3").Value = 1000
$ws.Range("F4").Value = 1.64
$ws.Range("G4").Value = 1.73
$ws.Range("H4").Value = 6.2
$ws.Range("I4").Value = 7.6
$ws.Range("J4").Value = 3.55
$ws.Range("K4").Value = 3.95
$ws.Range("P4").Value = 1.78
$ws.Range("R4").Value = 1.29
$ws.Range("T4").Value = 2.04
$ws.Range("U4").Value = 1.78
$ws.Range("V4").Value = 1.16
$ws.Range("W4").Value = 2.34
$ws.Range("Y4").Value = 980
$ws.Range("Z4").Value = 1000
$ws.Range("AF4").Value = 9.4
$ws.Range("AJ4").Value = 38
$ws.Range("AN4").Value = 15.5
$ws.Range("F5").Value = 1.69
$ws.Range("G5").Value = 1.71
$ws.Range("J5").Value = 4
$ws.Range("K5").Value = 4.1
$ws.Range("L5").Value = 1.41
$ws.Range("N5").Value = 3.95
$ws.Range("O5").Value = 1.32
$ws.Range("P5").Value = 1.98
$ws.Range("Q5").Value = 2
$ws.Range("R5").Value = 1.37
$ws.Range("S5").Value = 3.55
$ws.Range("T5").Value = 1.96
$ws.Range("U5").Value = 2
$ws.Range("X5").Value = 15
$ws.Range("Y5").Value = 20
$ws.Range("Z5").Value = 55
$ws.Range("AC5").Value = 8.800000000000001
$ws.Range("AF5").Value = 9.4
$ws.Range("AG5").Value = 10
$ws.Range("AH5").Value = 20
$ws.Range("AJ5").Value = 16
$ws.Range("AK5").Value = 17.5
$ws.Range("AO5").Value = 540
$ws.Range("F6").Value = 4.8
$ws.Range("G6").Value = 5.1
$ws.Range("H6").Value = 1.88
$ws.Range("I6").Value = 1.94
$ws.Range("J6").Value = 3.55
$ws.Range("K6").Value = 3.75
$ws.Range("L6").Value = 1.44
$ws.Range("S6").Value = 4
$ws.Range("T6").Value = 1.93
$ws.Range("W6").Value = 1.24
$ws.Range("X6").Value = 12
$ws.Range("Y6").Value = 8
$ws.Range("Z6").Value = 11
$ws.Range("AA6").Value = 21
$ws.Range("AB6").Value = 16
$ws.Range("AC6").Value = 8
$ws.Range("AD6").Value = 10.5
$ws.Range("AE6").Value = 22
$ws.Range("AF6").Value = 36
$ws.Range("AG6").Value = 21
$ws.Range("AH6").Value = 22
$ws.Range("AI6").Value = 42
$ws.Range("AJ6").Value = 120
$ws.Range("AK6").Value = 70
$ws.Range("AL6").Value = 85
$ws.Range("AM6").Value = 140
$ws.Range("AN6").Value = 90
$ws.Range("AO6").Value = 16
